{"js": "// Fix missing \"s\" on the academicBackground(s) loop tags: the opening\n// \"{-w:p academicBackground}\" and closing \"{/academicBackground}\" template\n// tags should read \"academicBackgrounds\" (plural) in both places.\nconst body = context.document.body;\n\nconst results = body.search(\"academicBackground\", {\n  matchCase: true,\n  matchWholeWord: false,\n});\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  // Append the missing \"s\" right after \"academicBackground\" without\n  // touching anything else around it (works for both the \"{-w:p\n  // academicBackground}\" opener and the \"{/academicBackground}\" closer).\n  results.items[i].insertText(\"s\", Word.InsertLocation.end);\n}\nawait context.sync();\n", "ps1": "# Fix missing \"s\" on the academicBackground(s) loop tags: the opening\n# \"{-w:p academicBackground}\" and closing \"{/academicBackground}\" template\n# tags should read \"academicBackgrounds\" (plural) in both places.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$find.Execute(\n    \"academicBackground\",   # FindText\n    $true,                  # MatchCase\n    $false,                 # MatchWholeWord\n    $false,                 # MatchWildcards\n    $false,                 # MatchSoundsLike\n    $false,                 # MatchAllWordForms\n    $true,                  # Forward\n    $wdFindContinue,        # Wrap\n    $false,                 # Format\n    \"academicBackgrounds\",  # ReplaceWith\n    $wdReplaceAll           # Replace\n) | Out-Null\n"}
